$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert two new rows before the existing row 2, shifting current data down
$ws.Rows.Item(2).Resize(2).Insert()
# Drop the formatting that Insert() inherited from the header row above
$ws.Rows.Item(2).Resize(2).ClearFormats()

# New student 1 (row 2)
$ws.Cells.Item(2,1).Value = 20330051920137
$ws.Cells.Item(2,2).Value = "MAYAHUA"
$ws.Cells.Item(2,3).Value = "XOCHIQUISQUI"
$ws.Cells.Item(2,4).Value = "DAMARIS"
$ws.Cells.Item(2,5).Value = "DISTINGUE LOS DIFERENTES TIPOS DE EMPRESA POR SU GIRO, ÁREAS FUNCIONALES, DOCUMENTACIÓN ADMINISTRATIVA Y RECURSOS"
$ws.Cells.Item(2,6).Value = "2ARHV"
$ws.Cells.Item(2,7).Value = 2

# New student 2 (row 3)
$ws.Cells.Item(3,1).Value = 19330051920130
$ws.Cells.Item(3,2).Value = "JIMENEZ"
$ws.Cells.Item(3,3).Value = "SANCHEZ"
$ws.Cells.Item(3,4).Value = "KAREN"
$ws.Cells.Item(3,5).Value = "EVALÚA EL DESEMPEÑO DE LA ORGANIZACIÓN UTILIZANDO HERRAMIENTAS DE CALIDAD"
$ws.Cells.Item(3,6).Value = "4ARHM"
$ws.Cells.Item(3,7).Value = 2

# The previously-existing rows (now rows 4-6) have their Reprobadas count reduced from 2 to 1
$ws.Cells.Item(4,7).Value = 1
$ws.Cells.Item(5,7).Value = 1
$ws.Cells.Item(6,7).Value = 1
